# Applies the "Updated cryptos list" GitHub Actions refresh:
# new Price (D) / Volume(1h) (E) figures, plus three rank swaps
# (rows 12/13 TRON<->WrappedEther, 19/21 Dai<->Uniswap via row 20,
# and 39/40 Aave<->MXToken) that keep column A (rank index) fixed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.751.28'
$ws.Range('E2').Value = '  -0.78%  '
$ws.Range('D3').Value = '1.928.45'
$ws.Range('E3').Value = '  -1.33%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.20'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9991'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.24%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4870'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2943'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06853'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.20'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.11%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '105.63'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.44%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07759'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.43%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.928.12'
$ws.Range('E13').Value = '  -1.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.333'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.96%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6981'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.37%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '274.26'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.94%  '
$ws.Range('D17').Value = '30.755.75'
$ws.Range('E17').Value = '  -0.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007693'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.23%  '
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9993'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('E20').Value = '  -1.25%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.595'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9999'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.481'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.844'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '164.39'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '19.47'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.49%  '
$ws.Range('E27').Value = '  -1.82%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1034'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.96%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.382'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.95%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.572'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.82%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.546'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.368'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.53%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04877'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7581'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.82%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.144'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.99%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9991'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.24%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.711'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.71%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01994'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.56%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.651'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.98%  '
$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '78.52'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.506'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.066'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.67%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8882'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4436'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.40%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '107.86'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.18%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.894'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.92%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9987'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.27%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '977.20'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.35%  '
$ws.Range('E49').Value = '  -1.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.14'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.42%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.182'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.54%  '
